$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 100007496
$ws.Cells.Item(18, 9).Value = 6992.8
$ws.Cells.Item(18, 11).Value = 6992.8
$ws.Cells.Item(18, 13).Value = -6708.8
$ws.Cells.Item(38, 8).Value = 6371.222
$ws.Cells.Item(38, 10).Value = 9483.166999999999
$ws.Cells.Item(38, 12).Value = 28449.501
$ws.Cells.Item(38, 14).Value = -29193.501
$ws.Cells.Item(42, 8).Value = 548.4286
$ws.Cells.Item(42, 9).Value = 600
$ws.Cells.Item(42, 11).Value = 1800
$ws.Cells.Item(42, 13).Value = -1570
$ws.Cells.Item(43, 8).Value = 5998.6665
$ws.Cells.Item(43, 10).Value = 5998
$ws.Cells.Item(43, 12).Value = 5998
$ws.Cells.Item(43, 14).Value = -6136
$ws.Cells.Item(51, 8).Value = 10803.2
$ws.Cells.Item(51, 10).Value = 10888.667
$ws.Cells.Item(51, 12).Value = 10888.667
$ws.Cells.Item(51, 14).Value = -11856.667
$ws.Cells.Item(53, 8).Value = 55556188
$ws.Cells.Item(53, 9).Value = 71429370
$ws.Cells.Item(53, 11).Value = 71429370
$ws.Cells.Item(53, 13).Value = -71428733
$ws.Cells.Item(64, 8).Value = 250008500
$ws.Cells.Item(64, 10).Value = 250008500
$ws.Cells.Item(64, 12).Value = 250008500
$ws.Cells.Item(64, 14).Value = -250008996
$ws.Cells.Item(67, 8).Value = 250008500
$ws.Cells.Item(67, 10).Value = 250008500
$ws.Cells.Item(67, 12).Value = 250008500
$ws.Cells.Item(67, 14).Value = -250010216
$ws.Cells.Item(86, 8).Value = 153849120
$ws.Cells.Item(86, 9).Value = 400001920
$ws.Cells.Item(86, 11).Value = 400001920
$ws.Cells.Item(86, 13).Value = -400000797
$ws.Cells.Item(89, 8).Value = 153849120
$ws.Cells.Item(89, 9).Value = 400001920
$ws.Cells.Item(89, 11).Value = 2000009600
$ws.Cells.Item(89, 13).Value = -2000003984
$ws.Cells.Item(123, 8).Value = 90998.336
$ws.Cells.Item(123, 10).Value = 90998.336
$ws.Cells.Item(123, 12).Value = 90998.336
$ws.Cells.Item(123, 14).Value = -100798.336
$ws.Cells.Item(126, 8).Value = 77738.664
$ws.Cells.Item(126, 10).Value = 77738.664
$ws.Cells.Item(126, 12).Value = 77738.664
$ws.Cells.Item(126, 14).Value = -87618.664
$ws.Cells.Item(129, 8).Value = 2497.125
$ws.Cells.Item(129, 9).Value = 1412.1428
$ws.Cells.Item(129, 10).Value = 2943.8823
$ws.Cells.Item(129, 11).Value = 4236.428400000001
$ws.Cells.Item(129, 12).Value = 8831.6469
$ws.Cells.Item(129, 13).Value = 763.5715999999993
$ws.Cells.Item(129, 14).Value = -18831.6469
$ws.Cells.Item(131, 8).Value = 527650.4
$ws.Cells.Item(131, 9).Value = 527650.4
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 1582951.2
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 13).Value = -1577911.2
$ws.Cells.Item(131, 14).Value = $null
$ws.Cells.Item(138, 8).Value = 3265.57
$ws.Cells.Item(138, 9).Value = 2775.5264
$ws.Cells.Item(138, 10).Value = 3380.5186
$ws.Cells.Item(138, 11).Value = 8326.5792
$ws.Cells.Item(138, 12).Value = 10141.5558
$ws.Cells.Item(138, 13).Value = -3186.5792
$ws.Cells.Item(138, 14).Value = -20421.5558
$ws.Cells.Item(141, 8).Value = 689.1429000000001
$ws.Cells.Item(141, 9).Value = 458
$ws.Cells.Item(141, 11).Value = 1374
$ws.Cells.Item(141, 13).Value = 3806

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(5, 8).Value = 195.8
$ws.Cells.Item(5, 9).Value = 182.25
$ws.Cells.Item(5, 11).Value = 182.25
$ws.Cells.Item(5, 13).Value = -70.25
$ws.Cells.Item(13, 8).Value = 9867.333000000001
$ws.Cells.Item(13, 9).Value = 9867.333000000001
$ws.Cells.Item(13, 11).Value = 9867.333000000001
$ws.Cells.Item(13, 13).Value = -9723.333000000001
$ws.Cells.Item(30, 8).Value = 2463
$ws.Cells.Item(30, 9).Value = 942.8570999999999
$ws.Cells.Item(30, 10).Value = 5123.25
$ws.Cells.Item(30, 11).Value = 942.8570999999999
$ws.Cells.Item(30, 12).Value = 5123.25
$ws.Cells.Item(30, 13).Value = -792.8570999999999
$ws.Cells.Item(30, 14).Value = -5423.25
$ws.Cells.Item(32, 8).Value = 5138.47
$ws.Cells.Item(32, 9).Value = 4415.0938
$ws.Cells.Item(32, 11).Value = 4415.0938
$ws.Cells.Item(32, 13).Value = -4128.0938
$ws.Cells.Item(61, 8).Value = 5245.154
$ws.Cells.Item(61, 9).Value = 4316.1177
$ws.Cells.Item(61, 10).Value = 7000
$ws.Cells.Item(61, 11).Value = 4316.1177
$ws.Cells.Item(61, 12).Value = 7000
$ws.Cells.Item(61, 13).Value = -4104.1177
$ws.Cells.Item(61, 14).Value = -7424
$ws.Cells.Item(74, 8).Value = 226076.2
$ws.Cells.Item(74, 9).Value = 467056.84
$ws.Cells.Item(74, 10).Value = 3632.5386
$ws.Cells.Item(74, 11).Value = 467056.84
$ws.Cells.Item(74, 12).Value = 3632.5386
$ws.Cells.Item(74, 13).Value = -466182.84
$ws.Cells.Item(74, 14).Value = -5380.5386
$ws.Cells.Item(77, 8).Value = 226076.2
$ws.Cells.Item(77, 9).Value = 467056.84
$ws.Cells.Item(77, 10).Value = 3632.5386
$ws.Cells.Item(77, 11).Value = 2335284.2
$ws.Cells.Item(77, 12).Value = 18162.693
$ws.Cells.Item(77, 13).Value = -2330916.2
$ws.Cells.Item(77, 14).Value = -26898.693
$ws.Cells.Item(102, 8).Value = 3714.7317
$ws.Cells.Item(102, 9).Value = 3625.8572
$ws.Cells.Item(102, 11).Value = 3625.8572
$ws.Cells.Item(102, 13).Value = -2003.8572
$ws.Cells.Item(122, 8).Value = 3298.0925
$ws.Cells.Item(122, 9).Value = 3190.5557
$ws.Cells.Item(122, 11).Value = 9571.667099999999
$ws.Cells.Item(122, 13).Value = -7121.667099999999
$ws.Cells.Item(130, 8).Value = 104999
$ws.Cells.Item(130, 10).Value = 104999
$ws.Cells.Item(130, 12).Value = 104999
$ws.Cells.Item(130, 14).Value = -115039
$ws.Cells.Item(131, 8).Value = 0
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 14).Value = $null
$ws.Cells.Item(136, 8).Value = 5245.154
$ws.Cells.Item(136, 9).Value = 4316.1177
$ws.Cells.Item(136, 10).Value = 7000
$ws.Cells.Item(136, 11).Value = 12948.3531
$ws.Cells.Item(136, 12).Value = 21000
$ws.Cells.Item(136, 13).Value = -10398.3531
$ws.Cells.Item(136, 14).Value = -26100

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 195.8
$ws.Cells.Item(4, 9).Value = 182.25
$ws.Cells.Item(4, 11).Value = 182.25
$ws.Cells.Item(4, 13).Value = -67.25
$ws.Cells.Item(99, 8).Value = 3682.5
$ws.Cells.Item(99, 9).Value = 2198.182
$ws.Cells.Item(99, 11).Value = 2198.182
$ws.Cells.Item(99, 13).Value = -700.1819999999998
$ws.Cells.Item(105, 8).Value = 20003628
$ws.Cells.Item(105, 9).Value = 1431981.8
$ws.Cells.Item(105, 10).Value = 41670548
$ws.Cells.Item(105, 11).Value = 1431981.8
$ws.Cells.Item(105, 12).Value = 41670548
$ws.Cells.Item(105, 13).Value = -1430234.8
$ws.Cells.Item(105, 14).Value = -41674042
$ws.Cells.Item(134, 8).Value = 3935.7144
$ws.Cells.Item(134, 9).Value = 3766.8333
$ws.Cells.Item(134, 10).Value = 4062.375
$ws.Cells.Item(134, 11).Value = 11300.4999
$ws.Cells.Item(134, 12).Value = 12187.125
$ws.Cells.Item(134, 13).Value = -8765.499899999999
$ws.Cells.Item(134, 14).Value = -17257.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 6069.9
$ws.Cells.Item(7, 9).Value = 6069.9
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 6069.9
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -5956.9
$ws.Cells.Item(7, 14).Value = $null
$ws.Cells.Item(31, 8).Value = 2757.39
$ws.Cells.Item(31, 9).Value = 2458.1428
$ws.Cells.Item(31, 11).Value = 2458.1428
$ws.Cells.Item(31, 13).Value = -2163.1428
$ws.Cells.Item(34, 8).Value = 2757.39
$ws.Cells.Item(34, 9).Value = 2458.1428
$ws.Cells.Item(34, 11).Value = 2458.1428
$ws.Cells.Item(34, 13).Value = -2256.1428
$ws.Cells.Item(41, 8).Value = 14177.083
$ws.Cells.Item(41, 9).Value = 10096.363
$ws.Cells.Item(41, 10).Value = 59065
$ws.Cells.Item(41, 11).Value = 10096.363
$ws.Cells.Item(41, 12).Value = 59065
$ws.Cells.Item(41, 13).Value = -9668.362999999999
$ws.Cells.Item(41, 14).Value = -59921
$ws.Cells.Item(50, 8).Value = 79000
$ws.Cells.Item(50, 10).Value = 79000
$ws.Cells.Item(50, 12).Value = 79000
$ws.Cells.Item(50, 14).Value = -80250
$ws.Cells.Item(51, 8).Value = 9999.5
$ws.Cells.Item(51, 9).Value = 9999.5
$ws.Cells.Item(51, 10).Value = 0
$ws.Cells.Item(51, 11).Value = 9999.5
$ws.Cells.Item(51, 12).Value = 0
$ws.Cells.Item(51, 13).Value = -9263.5
$ws.Cells.Item(51, 14).Value = $null
$ws.Cells.Item(60, 8).Value = 9138.375
$ws.Cells.Item(60, 9).Value = 6701.4
$ws.Cells.Item(60, 10).Value = 13200
$ws.Cells.Item(60, 11).Value = 6701.4
$ws.Cells.Item(60, 12).Value = 13200
$ws.Cells.Item(60, 13).Value = -6190.4
$ws.Cells.Item(60, 14).Value = -14222
$ws.Cells.Item(61, 8).Value = 9999.5
$ws.Cells.Item(61, 9).Value = 9999.5
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 9999.5
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -9651.5
$ws.Cells.Item(61, 14).Value = $null
$ws.Cells.Item(62, 8).Value = 5889107.5
$ws.Cells.Item(62, 9).Value = 11115077
$ws.Cells.Item(62, 11).Value = 11115077
$ws.Cells.Item(62, 13).Value = -11114453
$ws.Cells.Item(65, 8).Value = 5889107.5
$ws.Cells.Item(65, 9).Value = 11115077
$ws.Cells.Item(65, 11).Value = 55575385
$ws.Cells.Item(65, 13).Value = -55572265
$ws.Cells.Item(86, 8).Value = 3448
$ws.Cells.Item(86, 9).Value = 2518
$ws.Cells.Item(86, 10).Value = 4998
$ws.Cells.Item(86, 11).Value = 2518
$ws.Cells.Item(86, 12).Value = 4998
$ws.Cells.Item(86, 13).Value = -1395
$ws.Cells.Item(86, 14).Value = -7244
$ws.Cells.Item(89, 8).Value = 3448
$ws.Cells.Item(89, 9).Value = 2518
$ws.Cells.Item(89, 10).Value = 4998
$ws.Cells.Item(89, 11).Value = 12590
$ws.Cells.Item(89, 12).Value = 24990
$ws.Cells.Item(89, 13).Value = -6974
$ws.Cells.Item(89, 14).Value = -36222
$ws.Cells.Item(99, 8).Value = 8232.208000000001
$ws.Cells.Item(99, 9).Value = 10792
$ws.Cells.Item(99, 11).Value = 10792
$ws.Cells.Item(99, 13).Value = -9294
$ws.Cells.Item(107, 8).Value = 904.5294
$ws.Cells.Item(107, 9).Value = 918.26666
$ws.Cells.Item(107, 10).Value = 801.5
$ws.Cells.Item(107, 11).Value = 918.26666
$ws.Cells.Item(107, 12).Value = 801.5
$ws.Cells.Item(107, 13).Value = 1001.73334
$ws.Cells.Item(107, 14).Value = -4641.5
$ws.Cells.Item(126, 8).Value = 8232.208000000001
$ws.Cells.Item(126, 9).Value = 10792
$ws.Cells.Item(126, 11).Value = 32376
$ws.Cells.Item(126, 13).Value = -29906
$ws.Cells.Item(141, 8).Value = 730883
$ws.Cells.Item(141, 10).Value = 730883
$ws.Cells.Item(141, 12).Value = 730883
$ws.Cells.Item(141, 14).Value = -741243

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 1127.5625
$ws.Cells.Item(2, 9).Value = 925.0833
$ws.Cells.Item(2, 11).Value = 5550.4998
$ws.Cells.Item(2, 13).Value = -5437.4998
$ws.Cells.Item(5, 8).Value = 2051.577
$ws.Cells.Item(5, 9).Value = 428.125
$ws.Cells.Item(5, 11).Value = 1284.375
$ws.Cells.Item(5, 13).Value = -1172.375
$ws.Cells.Item(13, 8).Value = 99.5
$ws.Cells.Item(13, 9).Value = 99
$ws.Cells.Item(13, 11).Value = 297
$ws.Cells.Item(13, 13).Value = -129
$ws.Cells.Item(33, 8).Value = 82.88
$ws.Cells.Item(33, 9).Value = 76
$ws.Cells.Item(33, 10).Value = 85.052635
$ws.Cells.Item(33, 11).Value = 456
$ws.Cells.Item(33, 12).Value = 510.3158099999999
$ws.Cells.Item(33, 13).Value = -173
$ws.Cells.Item(33, 14).Value = -1076.31581
$ws.Cells.Item(113, 8).Value = 2443.3
$ws.Cells.Item(113, 10).Value = 2520.3333
$ws.Cells.Item(113, 12).Value = 7560.999899999999
$ws.Cells.Item(113, 14).Value = -11900.9999
$ws.Cells.Item(122, 8).Value = 1293.2941
$ws.Cells.Item(122, 10).Value = 1267.0667
$ws.Cells.Item(122, 12).Value = 11403.6003
$ws.Cells.Item(122, 14).Value = -16303.6003
$ws.Cells.Item(123, 8).Value = 1946.1538
$ws.Cells.Item(123, 9).Value = 750
$ws.Cells.Item(123, 11).Value = 2250
$ws.Cells.Item(123, 13).Value = 200
$ws.Cells.Item(135, 8).Value = 2051.577
$ws.Cells.Item(135, 9).Value = 428.125
$ws.Cells.Item(135, 11).Value = 3853.125
$ws.Cells.Item(135, 13).Value = -1318.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 611.4375
$ws.Cells.Item(2, 9).Value = 652.53845
$ws.Cells.Item(2, 10).Value = 433.33334
$ws.Cells.Item(2, 11).Value = 652.53845
$ws.Cells.Item(2, 12).Value = 433.33334
$ws.Cells.Item(2, 13).Value = -539.53845
$ws.Cells.Item(2, 14).Value = -659.33334
$ws.Cells.Item(14, 8).Value = 5000668.5
$ws.Cells.Item(14, 9).Value = 5000668.5
$ws.Cells.Item(14, 11).Value = 5000668.5
$ws.Cells.Item(14, 13).Value = -5000500.5
$ws.Cells.Item(62, 8).Value = 36415.832
$ws.Cells.Item(62, 9).Value = 24832.666
$ws.Cells.Item(62, 11).Value = 24832.666
$ws.Cells.Item(62, 13).Value = -24146.666
$ws.Cells.Item(65, 8).Value = 36415.832
$ws.Cells.Item(65, 9).Value = 24832.666
$ws.Cells.Item(65, 11).Value = 74497.99800000001
$ws.Cells.Item(65, 13).Value = -71065.99800000001
$ws.Cells.Item(70, 8).Value = 45642980
$ws.Cells.Item(70, 9).Value = 55783644
$ws.Cells.Item(70, 11).Value = 55783644
$ws.Cells.Item(70, 13).Value = -55783374
$ws.Cells.Item(73, 8).Value = 45642980
$ws.Cells.Item(73, 9).Value = 55783644
$ws.Cells.Item(73, 11).Value = 55783644
$ws.Cells.Item(73, 13).Value = -55782708
$ws.Cells.Item(80, 8).Value = 111114030
$ws.Cells.Item(80, 9).Value = 333335330
$ws.Cells.Item(80, 10).Value = 3381.5
$ws.Cells.Item(80, 11).Value = 333335330
$ws.Cells.Item(80, 12).Value = 3381.5
$ws.Cells.Item(80, 13).Value = -333334332
$ws.Cells.Item(80, 14).Value = -5377.5
$ws.Cells.Item(83, 8).Value = 111114030
$ws.Cells.Item(83, 9).Value = 333335330
$ws.Cells.Item(83, 10).Value = 3381.5
$ws.Cells.Item(83, 11).Value = 1666676650
$ws.Cells.Item(83, 12).Value = 16907.5
$ws.Cells.Item(83, 13).Value = -1666671658
$ws.Cells.Item(83, 14).Value = -26891.5
$ws.Cells.Item(102, 8).Value = 2296.9048
$ws.Cells.Item(102, 9).Value = 2354.3333
$ws.Cells.Item(102, 10).Value = 2220.3333
$ws.Cells.Item(102, 11).Value = 2354.3333
$ws.Cells.Item(102, 12).Value = 2220.3333
$ws.Cells.Item(102, 13).Value = -732.3332999999998
$ws.Cells.Item(102, 14).Value = -5464.3333
$ws.Cells.Item(103, 8).Value = 76998.60000000001
$ws.Cells.Item(103, 10).Value = 76998.60000000001
$ws.Cells.Item(103, 12).Value = 76998.60000000001
$ws.Cells.Item(103, 14).Value = -79342.60000000001
$ws.Cells.Item(136, 8).Value = 79999.5
$ws.Cells.Item(136, 10).Value = 79999.5
$ws.Cells.Item(136, 12).Value = 239998.5
$ws.Cells.Item(136, 14).Value = -245098.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 931.087
$ws.Cells.Item(16, 9).Value = 953.05554
$ws.Cells.Item(16, 11).Value = 953.05554
$ws.Cells.Item(16, 13).Value = -783.05554
$ws.Cells.Item(17, 9).Value = 20000
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 20000
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 13).Value = -19830
$ws.Cells.Item(17, 14).Value = $null
$ws.Cells.Item(22, 8).Value = 71430620
$ws.Cells.Item(22, 10).Value = 200002930
$ws.Cells.Item(22, 12).Value = 200002930
$ws.Cells.Item(22, 14).Value = -200003520
$ws.Cells.Item(27, 8).Value = 71430620
$ws.Cells.Item(27, 10).Value = 200002930
$ws.Cells.Item(27, 12).Value = 200002930
$ws.Cells.Item(27, 14).Value = -200003144
$ws.Cells.Item(46, 8).Value = 1474.7441
$ws.Cells.Item(46, 10).Value = 1361.6
$ws.Cells.Item(46, 12).Value = 1361.6
$ws.Cells.Item(46, 14).Value = -1737.6
$ws.Cells.Item(93, 8).Value = 2324
$ws.Cells.Item(93, 9).Value = 2155.8333
$ws.Cells.Item(93, 11).Value = 2155.8333
$ws.Cells.Item(93, 13).Value = -907.8332999999998
$ws.Cells.Item(100, 8).Value = 5346.722
$ws.Cells.Item(100, 9).Value = 4869.467
$ws.Cells.Item(100, 10).Value = 7733
$ws.Cells.Item(100, 11).Value = 4869.467
$ws.Cells.Item(100, 12).Value = 7733
$ws.Cells.Item(100, 13).Value = -4328.467
$ws.Cells.Item(100, 14).Value = -8815
$ws.Cells.Item(128, 8).Value = 48476
$ws.Cells.Item(128, 10).Value = 48476
$ws.Cells.Item(128, 12).Value = 48476
$ws.Cells.Item(128, 14).Value = -58436
$ws.Cells.Item(132, 8).Value = 3522.1365
$ws.Cells.Item(132, 9).Value = 2593.125
$ws.Cells.Item(132, 10).Value = 5999.5
$ws.Cells.Item(132, 11).Value = 7779.375
$ws.Cells.Item(132, 12).Value = 17998.5
$ws.Cells.Item(132, 13).Value = -5249.375
$ws.Cells.Item(132, 14).Value = -23058.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(9, 8).Value = 0
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 12).Value = 0
$ws.Cells.Item(9, 14).Value = $null
$ws.Cells.Item(54, 8).Value = 49487.5
$ws.Cells.Item(54, 10).Value = 49487.5
$ws.Cells.Item(54, 12).Value = 49487.5
$ws.Cells.Item(54, 14).Value = -50527.5
$ws.Cells.Item(62, 8).Value = 6569.4287
$ws.Cells.Item(62, 9).Value = 5331
$ws.Cells.Item(62, 10).Value = 14000
$ws.Cells.Item(62, 11).Value = 5331
$ws.Cells.Item(62, 12).Value = 14000
$ws.Cells.Item(62, 13).Value = -4707
$ws.Cells.Item(62, 14).Value = -15248
$ws.Cells.Item(65, 8).Value = 6569.4287
$ws.Cells.Item(65, 9).Value = 5331
$ws.Cells.Item(65, 10).Value = 14000
$ws.Cells.Item(65, 11).Value = 26655
$ws.Cells.Item(65, 12).Value = 70000
$ws.Cells.Item(65, 13).Value = -23535
$ws.Cells.Item(65, 14).Value = -76240
$ws.Cells.Item(81, 8).Value = 6783.25
$ws.Cells.Item(81, 10).Value = 6375
$ws.Cells.Item(81, 12).Value = 12750
$ws.Cells.Item(81, 14).Value = -14872
$ws.Cells.Item(84, 8).Value = 6783.25
$ws.Cells.Item(84, 10).Value = 6375
$ws.Cells.Item(84, 12).Value = 63750
$ws.Cells.Item(84, 14).Value = -74358
$ws.Cells.Item(96, 8).Value = 7246.25
$ws.Cells.Item(96, 9).Value = 6796.6
$ws.Cells.Item(96, 10).Value = 7995.6665
$ws.Cells.Item(96, 11).Value = 6796.6
$ws.Cells.Item(96, 12).Value = 7995.6665
$ws.Cells.Item(96, 13).Value = -5423.6
$ws.Cells.Item(96, 14).Value = -10741.6665
$ws.Cells.Item(107, 8).Value = 806.63635
$ws.Cells.Item(107, 10).Value = 1480
$ws.Cells.Item(107, 12).Value = 4440
$ws.Cells.Item(107, 14).Value = -8280
$ws.Cells.Item(113, 8).Value = 1018.8421
$ws.Cells.Item(113, 9).Value = 1180.3334
$ws.Cells.Item(113, 11).Value = 3541.0002
$ws.Cells.Item(113, 13).Value = -1371.0002
$ws.Cells.Item(122, 8).Value = 10002718
$ws.Cells.Item(122, 9).Value = 2137.7
$ws.Cells.Item(122, 11).Value = 6413.099999999999
$ws.Cells.Item(122, 13).Value = -3963.099999999999
$ws.Cells.Item(126, 8).Value = 2733.125
$ws.Cells.Item(126, 9).Value = 2695
$ws.Cells.Item(126, 11).Value = 8085
$ws.Cells.Item(126, 13).Value = -5615
$ws.Cells.Item(136, 8).Value = 27028652
$ws.Cells.Item(136, 9).Value = 31251382
$ws.Cells.Item(136, 10).Value = 3184.2
$ws.Cells.Item(136, 11).Value = 93754146
$ws.Cells.Item(136, 12).Value = 9552.599999999999
$ws.Cells.Item(136, 13).Value = -93751596
$ws.Cells.Item(136, 14).Value = -14652.6
